$wb = $excel.ActiveWorkbook

# Column F ("想去人数" - number of people wanting to attend) updates.
# These updates apply identically to sheet 1 ("展览") and sheet 4 ("全部类型"),
# which hold duplicated data, matching the source diff.
$updates = @{
    2  = 1106
    3  = 817
    4  = 273
    8  = 2041
    9  = 7598
    11 = 419
    12 = 351
    13 = 136
    14 = 402
    15 = 152
    16 = 7751
    18 = 1343
    19 = 151
    22 = 151
    23 = 309
    24 = 143
    25 = 163
    28 = 21
    29 = 411
    30 = 614
    31 = 53
    33 = 60
    34 = 79
    36 = 74
}

$sheetIndexes = @(1, 4)

foreach ($sheetIndex in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
